$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 13 with new sensor reading values
$ws.Range("A13").Value = 45874.45852196016
$ws.Range("D13").Value = 17.11
$ws.Range("E13").Value = 84.06
$ws.Range("F13").Value = 316.41
$ws.Range("G13").Value = 10.9
$ws.Range("J13").Value = "11:00:16"

# Delete the old row 14 entirely (shift cells up)
$ws.Rows("14").Delete()
